# Project_Grading_Rubric_Checklist.xlsx edit
# Adds a "Brainstorm Ideas" column (new column C) with per-row notes,
# updates the final-grade note text from 30/50 to 30/35, and adds a
# hyperlink on the Bitwise Operators brainstorm note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank column at C (old C/D shift right to D/E) ---
$ws.Columns("C").Insert()

# --- 2. Header cell for the new column ---
$ws.Range("C1").Value = "Brainstorm Ideas"

# --- 3. Brainstorm note text for each row ---
$notes = @{
    5  = "Include a lot of comments and be smart of how we name our variables"
    6  = "Include some equations from Aersp 309"
    8  = "Create a library for all our functions. Create header files for all our functions"
    9  = "Look for bugs (e.g. recognize if the user tries to input a negative amount of fuel)"
    11 = "Create a library for all our functions"
    14 = "Can represent multiple boolean states using a single bit- each bit represents a different variable. See this website."
    15 = "create if statements and for loops"
    16 = "Create a function"
    18 = "The rocket object accepts pointers of the sub-components. (e.g: using pointers, we can tell the rocket how many stages it has, what boosters it's using, etc.)"
    20 = "Create multiple functions that have different parameters."
    23 = "Have pre-programmed rocket parts that you can choose from"
    24 = "Create a class for the whole rocket as well as each of the individual components"
    25 = "Create objects from a class"
    27 = "Create a class that inherits properties of another one"
    28 = "Overriding inhereted functions"
    30 = "Option to save the output to a .txt file"
    32 = "Multithreading: Have one thread that constantly checks if there are enough equations and known variables to solve for some unknown variables."
}

foreach ($row in $notes.Keys) {
    $ws.Cells.Item($row, 3).Value = $notes[$row]
}

# --- 4. Format the "Brainstorm Ideas" block (C1:C32) ---
# Header cell: bold, filled, bordered (no wrap)
$ws.Range("C1").Interior.Color = 11911321
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").Borders.Weight = 2
$ws.Range("C1").Font.Bold = $true

# Column-title row cell: bold, filled, bordered, wrapped
$ws.Range("C2").Interior.Color = 11911321
$ws.Range("C2").Borders.LineStyle = 1
$ws.Range("C2").Borders.Weight = 2
$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").WrapText = $true

# Body cells C3:C32: filled, bordered, wrap text
$body = $ws.Range("C3:C32")
$body.Interior.Color = 11911321
$body.Borders.LineStyle = 1
$body.Borders.Weight = 2
$body.WrapText = $true

# C26 keeps the fill/border but no wrap (matches source formatting)
$ws.Range("C26").WrapText = $false

# C33/C34 (outside the box): plain wrap-text only, like the old D33 note style
$ws.Range("C33:C34").WrapText = $true

# --- 5. Hyperlink on the bitwise-operators brainstorm note ---
$ws.Hyperlinks.Add($ws.Range("C14"), "https://en.wikipedia.org/wiki/Bit_field", "", "", $notes[14])

# --- 6. Update the "full credit" note text (30/50 -> 30/35), now in column E ---
$ws.Range("E33").Value = "Get full credit (30% of the final grade) if you score more than or equal to 30/35"

# --- 6b. Stray total formula under the new (empty) Grade column ---
$ws.Range("D33").Formula = "=SUM(D3:D32)"

# --- 7. Column widths ---
$ws.Columns("C").ColumnWidth = 65.28515625
$ws.Columns("E").ColumnWidth = 168.28515625

Write-Output "done"
